$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Roll the five displayed fiscal periods forward by one year: drop the
# oldest columns data, shift D..G left from old E..H, and populate the
# new rightmost column H with the newly published 1401/12 period.

$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-03-12 (9)"
$ws.Range("E9").Value = "1400-02-29 (8)"
$ws.Range("F9").Value = "1401-02-27 (11)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "1402-02-28"

$ws.Range("D11").Value = 2082058
$ws.Range("E11").Value = 2021486
$ws.Range("F11").Value = 3291336
$ws.Range("G11").Value = 5091269
$ws.Range("H11").Value = 7339981

$ws.Range("D12").Value = -1348256
$ws.Range("E12").Value = -1476986
$ws.Range("F12").Value = -2129091
$ws.Range("G12").Value = -3513432
$ws.Range("H12").Value = -5901984

$ws.Range("D13").Value = 733802
$ws.Range("E13").Value = 544500
$ws.Range("F13").Value = 1162245
$ws.Range("G13").Value = 1577837
$ws.Range("H13").Value = 1437997

$ws.Range("D14").Value = -143429
$ws.Range("E14").Value = -176710
$ws.Range("F14").Value = -252233
$ws.Range("G14").Value = -436250
$ws.Range("H14").Value = -708803

$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

$ws.Range("D16").Value = -4177
$ws.Range("E16").Value = -16176
$ws.Range("F16").Value = -40380
$ws.Range("G16").Value = -46494
$ws.Range("H16").Value = -103603

$ws.Range("D17").Value = 586196
$ws.Range("E17").Value = 351614
$ws.Range("F17").Value = 869632
$ws.Range("G17").Value = 1095093
$ws.Range("H17").Value = 625591

$ws.Range("D18").Value = -241041
$ws.Range("E18").Value = -268487
$ws.Range("F18").Value = -324632
$ws.Range("G18").Value = -509952
$ws.Range("H18").Value = -693603

$ws.Range("D19").Value = 503854
$ws.Range("E19").Value = 1065102
$ws.Range("F19").Value = 1205651
$ws.Range("G19").Value = 2089882
$ws.Range("H19").Value = 3953286

$ws.Range("D20").Value = 849009
$ws.Range("E20").Value = 1148229
$ws.Range("F20").Value = 1750651
$ws.Range("G20").Value = 2675023
$ws.Range("H20").Value = 3885274

$ws.Range("D21").Value = -86152
$ws.Range("E21").Value = -47416
$ws.Range("F21").Value = -147685
$ws.Range("G21").Value = -141435
$ws.Range("H21").Value = 0

$ws.Range("D22").Value = 762857
$ws.Range("E22").Value = 1100813
$ws.Range("F22").Value = 1602966
$ws.Range("G22").Value = 2533588
$ws.Range("H22").Value = 3885274

$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 1117
$ws.Range("H23").Value = 1713

$ws.Range("D24").Value = 762857
$ws.Range("E24").Value = 1100813
$ws.Range("F24").Value = 1602966
$ws.Range("G24").Value = 2534705
$ws.Range("H24").Value = 3886987

$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = 971
$ws.Range("F25").Value = 707
$ws.Range("G25").Value = 1118
$ws.Range("H25").Value = 1714

$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 1134000
$ws.Range("F26").Value = 2268000
$ws.Range("G26").Value = 2268000
$ws.Range("H26").Value = 2268000

$ws.Range("D27").Value = 336
$ws.Range("E27").Value = 485
$ws.Range("F27").Value = 707
$ws.Range("G27").Value = 1118
$ws.Range("H27").Value = 1714
